# Regenerate orders with updated distance/sizes.
# Applies the token substitutions:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# across every text cell in the used range of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$startRow = $used.Row
$startCol = $used.Column

for ($r = $startRow; $r -lt ($startRow + $rowCount); $r++) {
    for ($c = $startCol; $c -lt ($startCol + $colCount); $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -is [string]) {
            $newVal = $val -replace "D64", "D69"
            $newVal = $newVal -replace "D51", "D55"
            $newVal = $newVal -replace "D80", "D86"
            $newVal = $newVal -replace "S30", "S31"
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
